$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet's tab title in workbook's <sheets> element and update header text
$ws.Name = "Through 2022-12-02"

# Update the "December (through 12-01)" label in A13 to "December (through 12-02)"
$ws.Range("A13").Value = "December (through 12-02)"

# Update December (row 13) monthly values
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 14
$ws.Range("I13").Value = 10

# Update Total (row 14) values
$ws.Range("C14").Value = 569
$ws.Range("D14").Value = 827
$ws.Range("E14").Value = 689
$ws.Range("F14").Value = 536
$ws.Range("G14").Value = 1272
$ws.Range("H14").Value = 1657
$ws.Range("I14").Value = 1526
